$wb = $excel.ActiveWorkbook

# --- Sheet "INFO" (sheet1.xml) ---
$wsInfo = $wb.Worksheets.Item("INFO")

# Swap row 3 and row 4 contents (Name:/Group: labels swap places)
$a3 = $wsInfo.Range("A3").Text
$b3 = $wsInfo.Range("B3").Text
$a4 = $wsInfo.Range("A4").Text
$b4 = $wsInfo.Range("B4").Text

$wsInfo.Range("A3").Value = $a4
$wsInfo.Range("B3").Value = $b4
$wsInfo.Range("A4").Value = $a3
$wsInfo.Range("B4").Value = $b3

# Update selection to K19
$wsInfo.Range("K19").Select() | Out-Null

# --- Sheet "SubCalc_template" (sheet2.xml) ---
$wsTemplate = $wb.Worksheets.Item("SubCalc_template")

# Swap A1 and B1 header contents (Name/Group columns swap places)
$a1 = $wsTemplate.Range("A1").Text
$b1 = $wsTemplate.Range("B1").Text

$wsTemplate.Range("A1").Value = $b1
$wsTemplate.Range("B1").Value = $a1

# Update selection to F29
$wsTemplate.Range("F29").Select() | Out-Null
